# Update column F (dSF) values for several rows per repulled data / mean calc fix.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value  = -6
$ws.Range("F8").Value  = -2
$ws.Range("F9").Value  = -5
$ws.Range("F10").Value = -5
$ws.Range("F11").Value = 4
$ws.Range("F13").Value = 3
$ws.Range("F14").Value = 0
$ws.Range("F17").Value = -3
